# Add two new columns, I ("I0") and J ("IF"), to the right of the
# existing H ("IP") column.
#   - I0 is a constant 1 for every data row.
#   - IF duplicates the value already present in column H for that row.
# Formatting of the new header cells (I1/J1) mirrors H1's style, and the
# sheet's used-range dimension grows from A1:H36 to A1:J36 automatically
# once the new cells are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine how many data rows are present (header in row 1, data from
# row 2 through the last used row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, "H").End(-4162).Row

# Header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (fill/font/border/alignment) onto the two new
# header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows: I is always 1, J mirrors column H.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, "I").Value = 1
    $ws.Cells.Item($r, "J").Value = $ws.Cells.Item($r, "H").Value2
}
